$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M17").ClearContents()
$ws.Range("H17").Value = 1400
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1400
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4200
$ws.Range("N17").Value = -4536

$ws.Range("H32").Value = 6473.125
$ws.Range("I32").Value = 5633.1665
$ws.Range("J32").Value = 8993
$ws.Range("K32").Value = 5633.1665
$ws.Range("L32").Value = 8993
$ws.Range("M32").Value = -5307.1665

$ws.Range("H40").Value = 1499.9166
$ws.Range("I40").Value = 1499.909
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1499.909
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1324.909

$ws.Range("H62").Value = 11696.5
$ws.Range("I62").Value = 14794.75
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 14794.75
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -14170.75
$ws.Range("N62").Value = -6748

$ws.Range("H65").Value = 11696.5
$ws.Range("I65").Value = 14794.75
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 73973.75
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -70853.75
$ws.Range("N65").Value = -33740

$ws.Range("N75").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0

$ws.Range("N78").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0

$ws.Range("H86").Value = 7897
$ws.Range("I86").Value = 8221.666999999999
$ws.Range("J86").Value = 7631.364
$ws.Range("K86").Value = 8221.666999999999
$ws.Range("L86").Value = 7631.364
$ws.Range("M86").Value = -7098.666999999999

$ws.Range("H89").Value = 7897
$ws.Range("I89").Value = 8221.666999999999
$ws.Range("J89").Value = 7631.364
$ws.Range("K89").Value = 41108.335
$ws.Range("L89").Value = 38156.82
$ws.Range("M89").Value = -35492.335

$ws.Range("H96").Value = 1142.1538
$ws.Range("I96").Value = 231.5
$ws.Range("J96").Value = 2599.2
$ws.Range("K96").Value = 694.5
$ws.Range("L96").Value = 7797.599999999999
$ws.Range("M96").Value = 678.5

$ws.Range("H138").Value = 4212.3105
$ws.Range("I138").Value = 4727.5
$ws.Range("J138").Value = 4129.88
$ws.Range("K138").Value = 14182.5
$ws.Range("L138").Value = 12389.64
$ws.Range("M138").Value = -9042.5
$ws.Range("N138").Value = -22669.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5915.5103
$ws.Range("I32").Value = 2466.5903
$ws.Range("J32").Value = 24999.533
$ws.Range("K32").Value = 2466.5903
$ws.Range("L32").Value = 24999.533
$ws.Range("M32").Value = -2179.5903
$ws.Range("N32").Value = -25573.533

$ws.Range("H45").Value = 1435.3636
$ws.Range("I45").Value = 1223.625
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1223.625
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -846.625

$ws.Range("H63").Value = 9954.75
$ws.Range("I63").Value = 9910
$ws.Range("J63").Value = 9999.5
$ws.Range("K63").Value = 9910
$ws.Range("L63").Value = 9999.5
$ws.Range("M63").Value = -9224

$ws.Range("H66").Value = 9954.75
$ws.Range("I66").Value = 9910
$ws.Range("J66").Value = 9999.5
$ws.Range("K66").Value = 49550
$ws.Range("L66").Value = 49997.5
$ws.Range("M66").Value = -46118

$ws.Range("N80").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0

$ws.Range("N83").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0

$ws.Range("H88").Value = 3326.5881
$ws.Range("I88").Value = 2948.1667
$ws.Range("J88").Value = 3533
$ws.Range("K88").Value = 2948.1667
$ws.Range("L88").Value = 3533
$ws.Range("M88").Value = -2542.1667

$ws.Range("H91").Value = 3326.5881
$ws.Range("I91").Value = 2948.1667
$ws.Range("J91").Value = 3533
$ws.Range("K91").Value = 2948.1667
$ws.Range("L91").Value = 3533
$ws.Range("M91").Value = -1544.1667

$ws.Range("H134").Value = 56383.332
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 56383.332
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 56383.332
$ws.Range("N134").Value = -66523.33199999999

$ws.Range("H135").Value = 63332.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 63332.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 63332.25
$ws.Range("N135").Value = -73472.25

$ws.Range("H137").Value = 89097
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 89097
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 89097
$ws.Range("N137").Value = -99297

$ws.Range("H140").Value = 39499.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39499.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39499.5
$ws.Range("N140").Value = -49859.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13665.308
$ws.Range("I82").Value = 7059
$ws.Range("J82").Value = 50000
$ws.Range("K82").Value = 7059
$ws.Range("L82").Value = 50000
$ws.Range("M82").Value = -6676
$ws.Range("N82").Value = -50766

$ws.Range("H85").Value = 13665.308
$ws.Range("I85").Value = 7059
$ws.Range("J85").Value = 50000
$ws.Range("K85").Value = 7059
$ws.Range("L85").Value = 50000
$ws.Range("M85").Value = -5733
$ws.Range("N85").Value = -52652

$ws.Range("H134").Value = 3975776.8
$ws.Range("I134").Value = 5266.875
$ws.Range("J134").Value = 9269790
$ws.Range("K134").Value = 15800.625
$ws.Range("L134").Value = 27809370
$ws.Range("M134").Value = -13265.625

$ws.Range("H135").Value = 99999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 99999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139

$ws.Range("H141").Value = 208332.67
$ws.Range("I141").Value = 45000
$ws.Range("J141").Value = 289999
$ws.Range("K141").Value = 45000
$ws.Range("L141").Value = 289999
$ws.Range("M141").Value = -39820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2139871.5
$ws.Range("I31").Value = 2779631.5
$ws.Range("J31").Value = 7337.8335
$ws.Range("K31").Value = 2779631.5
$ws.Range("L31").Value = 7337.8335
$ws.Range("M31").Value = -2779336.5

$ws.Range("H34").Value = 2139871.5
$ws.Range("I34").Value = 2779631.5
$ws.Range("J34").Value = 7337.8335
$ws.Range("K34").Value = 2779631.5
$ws.Range("L34").Value = 7337.8335
$ws.Range("M34").Value = -2779429.5

$ws.Range("H58").Value = 3798315.2
$ws.Range("I58").Value = 6560.778
$ws.Range("J58").Value = 6423376
$ws.Range("K58").Value = 6560.778
$ws.Range("L58").Value = 6423376
$ws.Range("M58").Value = -6357.778

$ws.Range("H132").Value = 7871.9585
$ws.Range("I132").Value = 3581.4
$ws.Range("J132").Value = 29324.75
$ws.Range("K132").Value = 10744.2
$ws.Range("L132").Value = 87974.25
$ws.Range("M132").Value = -8214.200000000001

$ws.Range("H136").Value = 3798315.2
$ws.Range("I136").Value = 6560.778
$ws.Range("J136").Value = 6423376
$ws.Range("K136").Value = 19682.334
$ws.Range("L136").Value = 19270128
$ws.Range("M136").Value = -17132.334

$ws.Range("H137").Value = 84466.664
$ws.Range("I137").Value = 60000
$ws.Range("J137").Value = 96700
$ws.Range("K137").Value = 60000
$ws.Range("L137").Value = 96700
$ws.Range("M137").Value = -54900
$ws.Range("N137").Value = -106900

$ws.Range("H140").Value = 79537.46000000001
$ws.Range("I140").Value = 70001
$ws.Range("J140").Value = 94795.8
$ws.Range("K140").Value = 70001
$ws.Range("L140").Value = 94795.8
$ws.Range("M140").Value = -64821

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 41872.168
$ws.Range("I2").Value = 77225.234
$ws.Range("J2").Value = 91.27273
$ws.Range("K2").Value = 77225.234
$ws.Range("L2").Value = 91.27273
$ws.Range("M2").Value = -77112.234

$ws.Range("H122").Value = 53193.81
$ws.Range("I122").Value = 72405.07000000001
$ws.Range("J122").Value = 5165.6665
$ws.Range("K122").Value = 217215.21
$ws.Range("L122").Value = 15496.9995
$ws.Range("M122").Value = -214765.21
$ws.Range("N122").Value = -20396.9995

$ws.Range("H126").Value = 14176.667
$ws.Range("I126").Value = 17877.273
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 53631.819
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -51161.819
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 20561.35
$ws.Range("I132").Value = 24775.455
$ws.Range("J132").Value = 15410.777
$ws.Range("K132").Value = 74326.36500000001
$ws.Range("L132").Value = 46232.331
$ws.Range("M132").Value = -71796.36500000001

$ws.Range("N133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws.Range("H135").Value = 87049.89999999999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 87049.89999999999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 87049.89999999999
$ws.Range("N135").Value = -97189.89999999999

$ws.Range("H139").Value = 180161.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 180161.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 180161.75
$ws.Range("N139").Value = -190441.75

$ws.Range("H141").Value = 123403.664
$ws.Range("I141").Value = 179999
$ws.Range("J141").Value = 95106
$ws.Range("K141").Value = 179999
$ws.Range("L141").Value = 95106
$ws.Range("M141").Value = -174819
$ws.Range("N141").Value = -105466

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8033.625
$ws.Range("I40").Value = 7752.7144
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 7752.7144
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -7616.7144

$ws.Range("H93").Value = 2573.7778
$ws.Range("I93").Value = 2114.3125
$ws.Range("J93").Value = 6249.5
$ws.Range("K93").Value = 2114.3125
$ws.Range("L93").Value = 6249.5
$ws.Range("M93").Value = -866.3125
$ws.Range("N93").Value = -8745.5

$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0

$ws.Range("H132").Value = 3341621.8
$ws.Range("I132").Value = 5845578.5
$ws.Range("J132").Value = 3012.7334
$ws.Range("K132").Value = 17536735.5
$ws.Range("L132").Value = 9038.200199999999
$ws.Range("M132").Value = -17534205.5
$ws.Range("N132").Value = -14098.2002

$ws.Range("H140").Value = 95737
$ws.Range("I140").Value = 45000
$ws.Range("J140").Value = 112649.336
$ws.Range("K140").Value = 45000
$ws.Range("L140").Value = 112649.336
$ws.Range("M140").Value = -39820
$ws.Range("N140").Value = -123009.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4506797.5
$ws.Range("I132").Value = 5378561.5
$ws.Range("J132").Value = 2683.3333
$ws.Range("K132").Value = 16135684.5
$ws.Range("L132").Value = 8049.999899999999
$ws.Range("M132").Value = -16133154.5

$ws.Range("H136").Value = 20290770
$ws.Range("I136").Value = 7247149
$ws.Range("J136").Value = 33334390
$ws.Range("K136").Value = 21741447
$ws.Range("L136").Value = 100003170
$ws.Range("M136").Value = -21738897
